$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28
$ws.Range("A28").Value = "Made map look better, changed a few of the mountains "

# Row 29
$ws.Range("A29").Value = 'Added a "stage"'

$ws.Range("B29").Formula = '="03/20/2017"'
$ws.Range("B29").Copy()
$ws.Range("B29").PasteSpecial(-4163)

$ws.Range("C29").Formula = '="03/25/2017"'
$ws.Range("C29").Copy()
$ws.Range("C29").PasteSpecial(-4163)

# Row 28 (E)
$ws.Range("E28").Value = "Simply trying to pretty up the map"

# Row 29 (E)
$ws.Range("E29").Value = "Added a labyrinth into the stage will add more"

# Row 30
$ws.Range("A30").Value = "Switch scene on objective completion"

$ws.Range("B30").Formula = '="03/27/2017"'
$ws.Range("B30").Copy()
$ws.Range("B30").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Row 28 date cells (B/C) - numeric, reuse the existing date style from B25/C25
$ws.Range("B25").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B28").Value = 42858

$ws.Range("C25").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 43011

$excel.CutCopyMode = $false

$ws.Range("F31").Select()
